## Prelim analysis with full Eurofins data
## Widen the data columns so the full values (client ID / sample type /
## collected-received-extracted-analyzed timestamps) are readable, mirroring
## the "best fit" column widths Excel applies when you double-click a column
## border (or select the columns and use Format > AutoFit Column Width).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B:C ("Client ID", "Sample Type") hold long text values
# (e.g. "Ongoing Precision & Recovery") -> best-fit width is 28 characters.
$ws.Range("B1:C1").EntireColumn.ColumnWidth = 27.166666666666668

# Columns D:F ("Collected", "Received", "Extracted") hold date/time stamps
# formatted as m/d/yy h:mm -> best-fit width ~13.86 characters.
$ws.Range("D1:F1").EntireColumn.ColumnWidth = 13

# Column G ("Analyzed") holds the same style of timestamp but with a couple
# of wider values -> best-fit width ~14.86 characters.
$ws.Range("G1:G1").EntireColumn.ColumnWidth = 14
